$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (SCD0224 -> SCD0015)
$ws.Name = "SCD0015"

# Update the TC_ID cell (DGS-239 -> SCD0015-002)
$ws.Range("B2").Value = "SCD0015-002"

# Widen column B to fit the new, longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 12.42578125

# Move the active selection from L3 to B3 (also resets the scrolled
# top-left cell back to the default)
$ws.Range("B3").Select() | Out-Null
